# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D: Price, E: Volume) keep their exact string
# representation instead of being auto-converted to numbers by Excel
# when values look numeric (e.g. "90.10", "1.00").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '35.371.24'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '1.853.21'
$ws.Range('E3').Value = '  +2.38%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '228.25'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  +2.48%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '42.95'
$ws.Range('E8').Value = '  +14.42%  '
$ws.Range('E9').Value = '  +5.80%  '
$ws.Range('D10').Value = '0.0692'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  +3.49%  '
$ws.Range('D12').Value = '2.122.07'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('D13').Value = '11.72'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').Value = '1.856.12'
$ws.Range('E14').Value = '  +2.46%  '
$ws.Range('D15').Value = '4.77'
$ws.Range('E15').Value = '  +7.51%  '
$ws.Range('D16').Value = '0.668'
$ws.Range('E16').Value = '  +5.31%  '
$ws.Range('D17').Value = '35.369.38'
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').Value = '70.07'
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('D19').Value = '246.85'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').Value = '12.27'
$ws.Range('E21').Value = '  +9.15%  '
$ws.Range('E22').Value = '  +15.35%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '172.03'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').Value = '7.95'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').Value = '17.91'
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '3.645.67'
$ws.Range('E29').Value = '  +50.05%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +7.72%  '
$ws.Range('D32').Value = '4.06'
$ws.Range('E32').Value = '  +3.03%  '
$ws.Range('D33').Value = '3.94'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('D34').Value = '0.0537'
$ws.Range('E34').Value = '  +3.09%  '
$ws.Range('D35').Value = '1.88'
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('D36').Value = '0.676'
$ws.Range('E36').Value = '  +3.47%  '
$ws.Range('D37').Value = '90.10'
$ws.Range('E37').Value = '  +11.31%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('E39').Value = '  +9.47%  '
$ws.Range('D40').Value = '1.341.89'
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('E41').Value = '  +4.08%  '
$ws.Range('D42').Value = '2.41'
$ws.Range('E42').Value = '  +2.27%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '1.27'
$ws.Range('E43').Value = '  +4.62%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '14.98'
$ws.Range('E44').Value = '  +8.08%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  +1.23%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '2.84'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('E47').Value = '  +3.76%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.019.65'
$ws.Range('E48').Value = '  +2.42%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '6.04'
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('D50').Value = '104.77'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.05%  '
